# Update the weekly units-completed report with the latest generation
# timestamp and the actual billed pricing figures now that the source data
# has been reconciled for this single work request / week-ending file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp (header blurb)
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# Total Billed Amount
$ws.Range("C8").Value = 15492.84

# Scope ID # value cleared (was "R25P7")
$ws.Range("G10").Value = ""

# Monday (07/14/2025) pricing
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55

# Tuesday (07/15/2025) pricing
$ws.Range("H22").Value = 478.55
$ws.Range("H23").Value = 478.55
$ws.Range("H24").Value = 476.4
$ws.Range("H25").Value = 1191
$ws.Range("H26").Value = 238.2
$ws.Range("H27").Value = 1548.3
$ws.Range("H28").Value = 1191
$ws.Range("H29").Value = 238.2
$ws.Range("H30").Value = 476.4
$ws.Range("H31").Value = 1429.2
$ws.Range("H32").Value = 7745.799999999999

# Wednesday (07/16/2025) pricing
$ws.Range("H37").Value = 478.55
$ws.Range("H38").Value = 63.44
$ws.Range("H39").Value = 282.51
$ws.Range("H40").Value = 55.18
$ws.Range("H41").Value = 350.53
$ws.Range("H42").Value = 1297.06
$ws.Range("H43").Value = 2527.27

# Thursday (07/17/2025) pricing
$ws.Range("H48").Value = 94.17
$ws.Range("H49").Value = 478.55
$ws.Range("H50").Value = 952.8
$ws.Range("H51").Value = 238.2
$ws.Range("H52").Value = 1763.72

# Friday (07/18/2025) pricing
$ws.Range("H57").Value = 476.4
$ws.Range("H58").Value = 476.4
$ws.Range("H59").Value = 476.4
$ws.Range("H60").Value = 476.4
$ws.Range("H61").Value = 476.4
$ws.Range("H62").Value = 476.4
$ws.Range("H63").Value = 119.1
$ws.Range("H64").Value = 2977.5
